# repull data, push all data, mean calculation
# Update column F (dSF) values for a set of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    8  = 2
    11 = 2
    14 = -1
    29 = 2
    36 = 2
    46 = 1
    56 = 0
    57 = 1
    59 = 0
    72 = 1
    80 = 4
    82 = 0
    85 = -1
    92 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
